$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dct:modified" timestamp in B21
$ws.Range("B21").Value = "2023-05-03T09:31:11+00:00"

# Append new vocabulary rows (27-33) pulled from the refreshed Google Sheet export
$ws.Range("A27").Value = "vocab:1003"
$ws.Range("B27").Value = "concentration of chemical substances in environment"
$ws.Range("D27").Value = "This is concentration level of a chemical in the environment"
$ws.Range("F27").Value = "subject"

$ws.Range("A28").Value = "vocab:1004"
$ws.Range("B28").Value = "concentration of alpha-HCH"
$ws.Range("F28").Value = "variable"

$ws.Range("A29").Value = "vocab:1005"
$ws.Range("B29").Value = "concentration of beta-HCH"
$ws.Range("F29").Value = "variable"

$ws.Range("A30").Value = "vocab:1006"
$ws.Range("B30").Value = "concentration of gamma-HCH"
$ws.Range("F30").Value = "variable"

$ws.Range("A31").Value = "vocab:1007"
$ws.Range("B31").Value = "concentration of delta-HCH"
$ws.Range("F31").Value = "variable"

$ws.Range("A32").Value = "vocab:1008"
$ws.Range("B32").Value = "concentration of sum of HCHs"
$ws.Range("F32").Value = "variable"

$ws.Range("A33").Value = "vocab:1009"
$ws.Range("B33").Value = "concentration of o,p'-DDE"
$ws.Range("F33").Value = "variable"
